$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet stores the "Price" (column D) and "Volume(1h)" (column E) figures
# as plain text (numeric-looking strings, e.g. "261.30", "0.68%"), not as
# numbers. A straight `.Value = "..."` assignment would make Excel infer a
# Number/Percentage type for these since they look numeric, which would
# change the stored cell type. To keep them as text we briefly force the
# Text number format before writing, then restore the cell's style to Normal
# so no formatting residue (beyond the refreshed values) is left behind.
$updates = @(
    @{ Cell = "D2";  Value = "261.29" },
    @{ Cell = "E2";  Value = "0.53%" },
    @{ Cell = "D3";  Value = "27.15" },
    @{ Cell = "E3";  Value = "1.05%" },
    @{ Cell = "D4";  Value = "4.706" },
    @{ Cell = "E4";  Value = "0.86%" },
    @{ Cell = "D5";  Value = "0.06213" },
    @{ Cell = "E5";  Value = "2.46%" },
    @{ Cell = "D6";  Value = "6.729" },
    @{ Cell = "D7";  Value = "0.8502" },
    @{ Cell = "E7";  Value = "-1.28%" },
    @{ Cell = "D8";  Value = "0.9073" },
    @{ Cell = "D9";  Value = "0.1402" },
    @{ Cell = "E9";  Value = "-0.20%" },
    @{ Cell = "D10"; Value = "0.04714" },
    @{ Cell = "E10"; Value = "-9.82%" },
    @{ Cell = "E11"; Value = "-0.84%" },
    @{ Cell = "D12"; Value = "0.03176" },
    @{ Cell = "E12"; Value = "1.85%" },
    @{ Cell = "D13"; Value = "0.09062" },
    @{ Cell = "E13"; Value = "-0.90%" },
    @{ Cell = "E14"; Value = "-0.02%" },
    @{ Cell = "D15"; Value = "0.0006179" },
    @{ Cell = "E15"; Value = "1.88%" },
    @{ Cell = "D16"; Value = "0.005991" },
    @{ Cell = "E16"; Value = "-2.15%" },
    @{ Cell = "D17"; Value = "3.466" },
    @{ Cell = "E17"; Value = "-0.55%" },
    @{ Cell = "D18"; Value = "3.173" },
    @{ Cell = "E18"; Value = "-0.20%" },
    @{ Cell = "D19"; Value = "2.178" },
    @{ Cell = "E19"; Value = "-0.34%" },
    @{ Cell = "E20"; Value = "-0.69%" },
    @{ Cell = "E21"; Value = "-0.59%" },
    @{ Cell = "D22"; Value = "4.091" },
    @{ Cell = "E22"; Value = "-0.18%" },
    @{ Cell = "D23"; Value = "0.04243" },
    @{ Cell = "E23"; Value = "0.20%" },
    @{ Cell = "D24"; Value = "0.001219" },
    @{ Cell = "E24"; Value = "0.15%" },
    @{ Cell = "D25"; Value = "0.004116" },
    @{ Cell = "E25"; Value = "1.92%" },
    @{ Cell = "E26"; Value = "0.10%" },
    @{ Cell = "D40"; Value = "0.03899" },
    @{ Cell = "E40"; Value = "0.50%" },
    @{ Cell = "D41"; Value = "0.1112" },
    @{ Cell = "E41"; Value = "-0.60%" },
    @{ Cell = "D42"; Value = "0.004134" },
    @{ Cell = "E42"; Value = "-0.36%" },
    @{ Cell = "E43"; Value = "-0.73%" },
    @{ Cell = "E44"; Value = "-9.29%" },
    @{ Cell = "E45"; Value = "-3.07%" },
    @{ Cell = "E46"; Value = "0.11%" },
    @{ Cell = "E47"; Value = "-35.77%" },
    @{ Cell = "D48"; Value = "0.1630" },
    @{ Cell = "E48"; Value = "20.44%" },
    @{ Cell = "E49"; Value = "0.11%" },
    @{ Cell = "E50"; Value = "0.11%" }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = "Normal"
}
